# Natmi following Dr Hou advice
# Recompute the LR-pair (Sema3a -> Plxna2) interaction table to include
# the "ECs" cluster as both a sending and a target cluster, expanding the
# 2x3 sending/target matrix to a full 3x3 (ECs, FAPs, sCs).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Sema3a"
$ws.Range("C2").Value = "Plxna2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.5674196666666667
$ws.Range("H2").Value = 1.702259
$ws.Range("I2").Value = 0.07864125446886469
$ws.Range("J2").Value = 0.07864125446886468
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 29.420614
$ws.Range("N2").Value = 88.261842
$ws.Range("O2").Value = 0.5865186809777162
$ws.Range("P2").Value = 0.5865186809777162
$ws.Range("Q2").Value = 16.69383498900867
$ws.Range("R2").Value = 150.244514901078
$ws.Range("S2").Value = 0.04612456484151145
$ws.Range("T2").Value = 0.04612456484151144

# Row 3: ECs -> FAPs
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Sema3a"
$ws.Range("C3").Value = "Plxna2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.5674196666666667
$ws.Range("H3").Value = 1.702259
$ws.Range("I3").Value = 0.07864125446886469
$ws.Range("J3").Value = 0.07864125446886468
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 4.080312
$ws.Range("N3").Value = 12.240936
$ws.Range("O3").Value = 0.08134361887272465
$ws.Range("P3").Value = 0.08134361887272466
$ws.Range("Q3").Value = 2.315249274936
$ws.Range("R3").Value = 20.837243474424
$ws.Range("S3").Value = 0.006396964231188284
$ws.Range("T3").Value = 0.006396964231188284

# Row 4: ECs -> sCs
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Sema3a"
$ws.Range("C4").Value = "Plxna2"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.5674196666666667
$ws.Range("H4").Value = 1.702259
$ws.Range("I4").Value = 0.07864125446886469
$ws.Range("J4").Value = 0.07864125446886468
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 16.660501
$ws.Range("N4").Value = 49.981503
$ws.Range("O4").Value = 0.3321377001495591
$ws.Range("P4").Value = 0.3321377001495591
$ws.Range("Q4").Value = 9.453495923919666
$ws.Range("R4").Value = 85.081463315277
$ws.Range("S4").Value = 0.02611972539616496
$ws.Range("T4").Value = 0.02611972539616495

# Row 5: FAPs -> ECs
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Sema3a"
$ws.Range("C5").Value = "Plxna2"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.7227763333333334
$ws.Range("H5").Value = 2.168329
$ws.Range("I5").Value = 0.1001728366019618
$ws.Range("J5").Value = 0.1001728366019618
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 29.420614
$ws.Range("N5").Value = 88.261842
$ws.Range("O5").Value = 0.5865186809777162
$ws.Range("P5").Value = 0.5865186809777162
$ws.Range("Q5").Value = 21.26452351133533
$ws.Range("R5").Value = 191.380711602018
$ws.Range("S5").Value = 0.05875323999357893
$ws.Range("T5").Value = 0.05875323999357893

# Row 6: FAPs -> FAPs
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Sema3a"
$ws.Range("C6").Value = "Plxna2"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.7227763333333334
$ws.Range("H6").Value = 2.168329
$ws.Range("I6").Value = 0.1001728366019618
$ws.Range("J6").Value = 0.1001728366019618
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 4.080312
$ws.Range("N6").Value = 12.240936
$ws.Range("O6").Value = 0.08134361887272465
$ws.Range("P6").Value = 0.08134361887272466
$ws.Range("Q6").Value = 2.949152946216
$ws.Range("R6").Value = 26.542376515944
$ws.Range("S6").Value = 0.008148421041949702
$ws.Range("T6").Value = 0.008148421041949704

# Row 7: FAPs -> sCs
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Sema3a"
$ws.Range("C7").Value = "Plxna2"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.7227763333333334
$ws.Range("H7").Value = 2.168329
$ws.Range("I7").Value = 0.1001728366019618
$ws.Range("J7").Value = 0.1001728366019618
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 16.660501
$ws.Range("N7").Value = 49.981503
$ws.Range("O7").Value = 0.3321377001495591
$ws.Range("P7").Value = 0.3321377001495591
$ws.Range("Q7").Value = 12.04181582427633
$ws.Range("R7").Value = 108.376342418487
$ws.Range("S7").Value = 0.03327117556643317
$ws.Range("T7").Value = 0.03327117556643317

# Row 8: sCs -> ECs
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Sema3a"
$ws.Range("C8").Value = "Plxna2"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 5.925096666666668
$ws.Range("H8").Value = 17.77529
$ws.Range("I8").Value = 0.8211859089291735
$ws.Range("J8").Value = 0.8211859089291734
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 29.420614
$ws.Range("N8").Value = 88.261842
$ws.Range("O8").Value = 0.5865186809777162
$ws.Range("P8").Value = 0.5865186809777162
$ws.Range("Q8").Value = 174.3199819426867
$ws.Range("R8").Value = 1568.87983748418
$ws.Range("S8").Value = 0.4816408761426259
$ws.Range("T8").Value = 0.4816408761426258

# Row 9: sCs -> FAPs
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Sema3a"
$ws.Range("C9").Value = "Plxna2"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 5.925096666666668
$ws.Range("H9").Value = 17.77529
$ws.Range("I9").Value = 0.8211859089291735
$ws.Range("J9").Value = 0.8211859089291734
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 4.080312
$ws.Range("N9").Value = 12.240936
$ws.Range("O9").Value = 0.08134361887272465
$ws.Range("P9").Value = 0.08134361887272466
$ws.Range("Q9").Value = 24.17624303016001
$ws.Range("R9").Value = 217.5861872714401
$ws.Range("S9").Value = 0.06679823359958667
$ws.Range("T9").Value = 0.06679823359958667

# Row 10: sCs -> sCs
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Sema3a"
$ws.Range("C10").Value = "Plxna2"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 5.925096666666668
$ws.Range("H10").Value = 17.77529
$ws.Range("I10").Value = 0.8211859089291735
$ws.Range("J10").Value = 0.8211859089291734
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 16.660501
$ws.Range("N10").Value = 49.981503
$ws.Range("O10").Value = 0.3321377001495591
$ws.Range("P10").Value = 0.3321377001495591
$ws.Range("Q10").Value = 98.71507894009669
$ws.Range("R10").Value = 888.4357104608702
$ws.Range("S10").Value = 0.272746799186961
$ws.Range("T10").Value = 0.272746799186961
